$d = $word.ActiveDocument

# Locate the run ", 3000-4000 words" that follows
# "Critical summary of previous research and show why my dissertation is needed"
$range = $d.Content
$found = $range.Find.Execute(", 3000-4000 words", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the text ', 3000-4000 words' in the document."
}

$start = $range.Start

# Update the word-count figures in place (same overall length, so character
# offsets below stay valid): "3000" -> "4000" and "4000" -> "4500".
$range.Text = ", 4000-4500 words"

# The single run that used to hold ", 3000-4000 words" has now been merged
# with its neighbouring runs ("...is needed" before it and ")" after it)
# because they all share identical (empty) run formatting. Re-split the text
# into three separate runs, matching the target layout:
#   ", 4000"  |  "-4500"  |  " words"
# by toggling a character formatting property off/on across the two outer
# segments - this forces new run boundaries at start, start+6, start+11 and
# start+17 while leaving the middle segment ("-4500") as a clean, separate
# run with no explicit run properties.
$segA = $d.Range($start, $start + 6)          # ", 4000"
$segA.Font.Bold = 1
$segA.Font.Bold = 0

$segC = $d.Range($start + 11, $start + 17)    # " words"
$segC.Font.Bold = 1
$segC.Font.Bold = 0
